# Commit: "added new project edtach"
# Rewrites the bullet text on the "Product Insights & Impact" slide (6) and
# the "Store & Location Insights" slide (7): the literal "- " / "->" style
# list-prefix characters that had been typed by hand are dropped (the
# placeholder's own bullet formatting already supplies them), and the text
# is retyped/cleaned up.

$p = $ppt.ActivePresentation

# --- Slide 6: "Product Insights & Impact" ---------------------------------
$s6 = $p.Slides.Item(6)
$body6 = $s6.Shapes.Item(2).TextFrame.TextRange

# Paragraph 1: "Top sellers: Lego Bricks, Colorbuds" - keep "Colorbuds" as
# its own run (it's the part flagged by the proofing tools in the source deck).
$para6_1 = $body6.Paragraphs(1)
$para6_1.Text = "Top sellers: Lego Bricks, Colorbuds"
$run6_1b = $para6_1.Characters(27, 9)
$run6_1b.Text = "Colorbuds"

# Paragraph 2: "Colorbuds: most profitable product (2022-2023)"
$para6_2 = $body6.Paragraphs(2)
$para6_2.Text = "Colorbuds: most profitable product (2022-2023)"
$run6_2a = $para6_2.Characters(1, 9)
$run6_2a.Text = "Colorbuds"

# Paragraph 3 & 4: simple re-typed text, arrow glyph dropped
$body6.Paragraphs(3).Text = "Insights enabled targeted campaigns"
$body6.Paragraphs(4).Text = "Result: 12% increase in marketing ROI"

# --- Slide 7: "Store & Location Insights" ----------------------------------
$s7 = $p.Slides.Item(7)
$body7 = $s7.Shapes.Item(2).TextFrame.TextRange

$body7.Paragraphs(1).Text = "Identified top and least performing stores"
$body7.Paragraphs(2).Text = "Tracked performance trends (YoY, quarterly)"
$body7.Paragraphs(3).Text = "Prioritized high-revenue cities (e.g., Ciudad de Mexico)"
$body7.Paragraphs(4).Text = "Result: 8% increase in total revenue"
